$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 89999
$ws.Range("J3").Value = 89999
$ws.Range("L3").Value = 89999
$ws.Range("N3").Value = -90227

$ws.Range("H51").Value = 127719720
$ws.Range("I51").Value = 340580540
$ws.Range("J51").Value = 3230
$ws.Range("K51").Value = 340580540
$ws.Range("L51").Value = 3230
$ws.Range("M51").Value = -340580056
$ws.Range("N51").Value = -4198

$ws.Range("H98").Value = 1401.25
$ws.Range("I98").Value = 1453.8096
$ws.Range("K98").Value = 1453.8096
$ws.Range("M98").Value = 44.19039999999995

$ws.Range("H102").Value = 89999
$ws.Range("J102").Value = 89999
$ws.Range("L102").Value = 89999
$ws.Range("N102").Value = -96489

$ws.Range("H122").Value = 1401.25
$ws.Range("I122").Value = 1453.8096
$ws.Range("K122").Value = 4361.4288
$ws.Range("M122").Value = -1911.4288

$ws.Range("H132").Value = 20002538
$ws.Range("I132").Value = 24393126
$ws.Range("K132").Value = 73179378
$ws.Range("M132").Value = -73176848

$ws.Range("H135").Value = 2063.081
$ws.Range("I135").Value = 960.4231
$ws.Range("K135").Value = 8643.8079
$ws.Range("M135").Value = -6108.8079

$ws.Range("H137").Value = 5765.6855
$ws.Range("J137").Value = 9799.714
$ws.Range("L137").Value = 29399.142
$ws.Range("N137").Value = -34499.142

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1668.4736
$ws.Range("J2").Value = 1132.2222
$ws.Range("L2").Value = 1132.2222
$ws.Range("N2").Value = -1358.2222

$ws.Range("H5").Value = 3831.4285
$ws.Range("I5").Value = 4871.1816
$ws.Range("J5").Value = 19
$ws.Range("K5").Value = 4871.1816
$ws.Range("L5").Value = 19
$ws.Range("M5").Value = -4759.1816
$ws.Range("N5").Value = -243

$ws.Range("H32").Value = 6775.9062
$ws.Range("I32").Value = 6126.3076
$ws.Range("K32").Value = 6126.3076
$ws.Range("M32").Value = -5839.3076

$ws.Range("H61").Value = 10973.257
$ws.Range("I61").Value = 7018.933
$ws.Range("K61").Value = 7018.933
$ws.Range("M61").Value = -6806.933

$ws.Range("H116").Value = 1668.4736
$ws.Range("J116").Value = 1132.2222
$ws.Range("L116").Value = 1132.2222
$ws.Range("N116").Value = -5720.2222

$ws.Range("H122").Value = 4523.8237
$ws.Range("I122").Value = 4429.1787
$ws.Range("K122").Value = 13287.5361
$ws.Range("M122").Value = -10837.5361

$ws.Range("H132").Value = 1826.6786
$ws.Range("I132").Value = 1746.02
$ws.Range("K132").Value = 5238.059999999999
$ws.Range("M132").Value = -2708.059999999999

$ws.Range("H136").Value = 10973.257
$ws.Range("I136").Value = 7018.933
$ws.Range("K136").Value = 21056.799
$ws.Range("M136").Value = -18506.799

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1668.4736
$ws.Range("J3").Value = 1132.2222
$ws.Range("L3").Value = 1132.2222
$ws.Range("N3").Value = -1360.2222

$ws.Range("H4").Value = 3831.4285
$ws.Range("I4").Value = 4871.1816
$ws.Range("J4").Value = 19
$ws.Range("K4").Value = 4871.1816
$ws.Range("L4").Value = 19
$ws.Range("M4").Value = -4756.1816
$ws.Range("N4").Value = -249

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2550.1765
$ws.Range("I31").Value = 2256.8
$ws.Range("J31").Value = 4750.5
$ws.Range("K31").Value = 2256.8
$ws.Range("L31").Value = 4750.5
$ws.Range("M31").Value = -1961.8
$ws.Range("N31").Value = -5340.5

$ws.Range("H34").Value = 2550.1765
$ws.Range("I34").Value = 2256.8
$ws.Range("J34").Value = 4750.5
$ws.Range("K34").Value = 2256.8
$ws.Range("L34").Value = 4750.5
$ws.Range("M34").Value = -2054.8
$ws.Range("N34").Value = -5154.5

$ws.Range("H99").Value = 6857
$ws.Range("I99").Value = 6633.222
$ws.Range("K99").Value = 6633.222
$ws.Range("M99").Value = -5135.222

$ws.Range("H126").Value = 6857
$ws.Range("I126").Value = 6633.222
$ws.Range("K126").Value = 19899.666
$ws.Range("M126").Value = -17429.666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 1787.125
$ws.Range("J103").Value = 4250
$ws.Range("L103").Value = 12750
$ws.Range("N103").Value = -14508

$ws.Range("H115").Value = 8762.5
$ws.Range("I115").Value = 10000
$ws.Range("J115").Value = 5050
$ws.Range("K115").Value = 30000
$ws.Range("L115").Value = 15150
$ws.Range("M115").Value = -28825
$ws.Range("N115").Value = -17500

$ws.Range("H129").Value = 91739.55
$ws.Range("I129").Value = 100474
$ws.Range("J129").Value = 4395
$ws.Range("K129").Value = 301422
$ws.Range("L129").Value = 13185
$ws.Range("M129").Value = -296422
$ws.Range("N129").Value = -23185

$ws.Range("H131").Value = 171138.67
$ws.Range("I131").Value = 1429542.9
$ws.Range("J131").Value = 1738.1154
$ws.Range("K131").Value = 4288628.699999999
$ws.Range("L131").Value = 5214.3462
$ws.Range("M131").Value = -4283588.699999999
$ws.Range("N131").Value = -15294.3462

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 18543.566
$ws.Range("I102").Value = 1697.48
$ws.Range("K102").Value = 1697.48
$ws.Range("M102").Value = -75.48000000000002

$ws.Range("H126").Value = 4527.75
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()

$ws.Range("H132").Value = 7410683
$ws.Range("I132").Value = 8336621.5
$ws.Range("K132").Value = 25009864.5
$ws.Range("M132").Value = -25007334.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2864.7097
$ws.Range("I93").Value = 3054.9583
$ws.Range("K93").Value = 3054.9583
$ws.Range("M93").Value = -1806.9583

$ws.Range("H122").Value = 5618.0347
$ws.Range("J122").Value = 5677.2
$ws.Range("L122").Value = 17031.6
$ws.Range("N122").Value = -21931.6

$ws.Range("H132").Value = 3879.9722
$ws.Range("I132").Value = 3737.111
$ws.Range("J132").Value = 4022.8333
$ws.Range("K132").Value = 11211.333
$ws.Range("L132").Value = 12068.4999
$ws.Range("M132").Value = -8681.332999999999
$ws.Range("N132").Value = -17128.4999

$ws.Range("H136").Value = 5196.357
$ws.Range("I136").Value = 4419.92
$ws.Range("K136").Value = 13259.76
$ws.Range("M136").Value = -10709.76

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 8169
$ws.Range("I4").Value = 1915
$ws.Range("J4").Value = 17550
$ws.Range("K4").Value = 1915
$ws.Range("L4").Value = 17550
$ws.Range("M4").Value = -1802
$ws.Range("N4").Value = -17776

$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()

$ws.Range("H126").Value = 3755.2727
$ws.Range("J126").Value = 3895.8
$ws.Range("L126").Value = 11687.4
$ws.Range("N126").Value = -16627.4
